$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.2210183333333333
$ws.Range("H2").Value = 0.6630550000000001
$ws.Range("I2").Value = 0.111623447733668
$ws.Range("J2").Value = 0.111623447733668
$ws.Range("M2").Value = 0.007742333333333334
$ws.Range("N2").Value = 0.023227
$ws.Range("O2").Value = 0.001217676423630818
$ws.Range("P2").Value = 0.001217676423630818
$ws.Range("Q2").Value = 0.001711197609444445
$ws.Range("R2").Value = 0.015400778485
$ws.Range("S2").Value = 0.0001359212406296744
$ws.Range("T2").Value = 0.0001359212406296744
# Row 3
$ws.Range("G3").Value = 0.2210183333333333
$ws.Range("H3").Value = 0.6630550000000001
$ws.Range("I3").Value = 0.111623447733668
$ws.Range("J3").Value = 0.111623447733668
$ws.Range("O3").Value = 0.9497929577862038
$ws.Range("P3").Value = 0.9497929577862039
$ws.Range("Q3").Value = 1.33474164999
$ws.Range("R3").Value = 12.01267484991
$ws.Range("S3").Value = 0.1060191645812543
$ws.Range("T3").Value = 0.1060191645812543
# Row 4
$ws.Range("G4").Value = 0.2210183333333333
$ws.Range("H4").Value = 0.6630550000000001
$ws.Range("I4").Value = 0.111623447733668
$ws.Range("J4").Value = 0.111623447733668
$ws.Range("M4").Value = 0.3114883333333334
$ws.Range("N4").Value = 0.9344650000000001
$ws.Range("O4").Value = 0.04898936579016542
$ws.Range("P4").Value = 0.04898936579016543
$ws.Range("Q4").Value = 0.06884463228611112
$ws.Range("R4").Value = 0.6196016905750001
$ws.Range("S4").Value = 0.005468361911784073
$ws.Range("T4").Value = 0.005468361911784074
# Row 5
$ws.Range("I5").Value = 0.5186760166697389
$ws.Range("J5").Value = 0.5186760166697389
$ws.Range("M5").Value = 0.007742333333333334
$ws.Range("N5").Value = 0.023227
$ws.Range("O5").Value = 0.001217676423630818
$ws.Range("P5").Value = 0.001217676423630818
$ws.Range("Q5").Value = 0.007951350525555555
$ws.Range("R5").Value = 0.07156215473000001
$ws.Range("S5").Value = 0.0006315795570014862
$ws.Range("T5").Value = 0.0006315795570014862
# Row 6
$ws.Range("I6").Value = 0.5186760166697389
$ws.Range("J6").Value = 0.5186760166697389
$ws.Range("O6").Value = 0.9497929577862038
$ws.Range("P6").Value = 0.9497929577862039
$ws.Range("S6").Value = 0.4926348280055177
$ws.Range("T6").Value = 0.4926348280055177
# Row 7
$ws.Range("I7").Value = 0.5186760166697389
$ws.Range("J7").Value = 0.5186760166697389
$ws.Range("M7").Value = 0.3114883333333334
$ws.Range("N7").Value = 0.9344650000000001
$ws.Range("O7").Value = 0.04898936579016542
$ws.Range("P7").Value = 0.04898936579016543
$ws.Range("Q7").Value = 0.3198974800388889
$ws.Range("S7").Value = 0.02540960910721978
$ws.Range("T7").Value = 0.02540960910721978
# Row 8
$ws.Range("I8").Value = 0.369700535596593
$ws.Range("J8").Value = 0.369700535596593
$ws.Range("M8").Value = 0.007742333333333334
$ws.Range("N8").Value = 0.023227
$ws.Range("O8").Value = 0.001217676423630818
$ws.Range("P8").Value = 0.001217676423630818
$ws.Range("Q8").Value = 0.005667542846666667
$ws.Range("R8").Value = 0.05100788562000001
$ws.Range("S8").Value = 0.0004501756259996573
$ws.Range("T8").Value = 0.0004501756259996573
# Row 9
$ws.Range("I9").Value = 0.369700535596593
$ws.Range("J9").Value = 0.369700535596593
$ws.Range("O9").Value = 0.9497929577862038
$ws.Range("P9").Value = 0.9497929577862039
$ws.Range("S9").Value = 0.3511389651994318
$ws.Range("T9").Value = 0.3511389651994318
# Row 10
$ws.Range("I10").Value = 0.369700535596593
$ws.Range("J10").Value = 0.369700535596593
$ws.Range("M10").Value = 0.3114883333333334
$ws.Range("N10").Value = 0.9344650000000001
$ws.Range("O10").Value = 0.04898936579016542
$ws.Range("P10").Value = 0.04898936579016543
$ws.Range("Q10").Value = 0.2280156897666667
$ws.Range("R10").Value = 2.052141207900001
$ws.Range("S10").Value = 0.01811139477116156
$ws.Range("T10").Value = 0.01811139477116157
